# Tijdschrijfformulier update - KBSa les 20/11/2023
# Adds a "KBS a les" time-tracking row (20-11-2023, 120 minuten, "Les") for
# every student who attended that class: Marvin, Luuk and Jochem.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Add-LesRow {
    param($ws, $newRow, $sourceRow)

    $srcA = "A" + $sourceRow
    $srcD = "D" + $sourceRow
    $dstA = "A" + $newRow
    $dstB = "B" + $newRow
    $dstC = "C" + $newRow
    $dstD = "D" + $newRow

    # Write the values first so the existing SUM(C10:C152) formulas pick the
    # new minutes up when the workbook recalculates.
    $ws.Range($dstA).Value = $ws.Range($srcA).Value2
    $ws.Range($dstB).Value = 45250
    $ws.Range($dstC).Value = 120
    $ws.Range($dstD).Value = $ws.Range($srcD).Value2

    # Then copy the formatting (incl. the date number format on column B)
    # from an existing "KBS a les" row so no new cell styles are created.
    $srcRange = "A" + $sourceRow + ":D" + $sourceRow
    $dstRange = "A" + $newRow + ":D" + $newRow
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial($xlPasteFormats)
}

$wsMarvin = $wb.Worksheets.Item("Marvin")
$wsLuuk = $wb.Worksheets.Item("Luuk")
$wsJochem = $wb.Worksheets.Item("Jochem")

Add-LesRow $wsMarvin 16 15
Add-LesRow $wsLuuk 16 12
Add-LesRow $wsJochem 19 12

# Update the active sheet/selection state. Luuk was the previously active
# tab; the new active tab is Marvin, with Luuk and Jochem's selections
# moved onto their freshly added rows.
$wsLuuk.Range("A16:D16").Select()
$wsJochem.Range("A19:D19").Select()
$wsMarvin.Range("C19").Select()
